# Fix the typo "Fenrnao" -> "Fenrnando" in the author's name on the
# cover line ("Andres Fenrnao Galvis" -> "Andres Fenrnando Galvis").
#
# The canonical OOXML diff shows this text being produced by three
# separate runs ("Andres Fenrna" / "nd" / "o Galvis") with the
# document's single "_GoBack" bookmark sitting right between the
# inserted "nd" run and the following "o Galvis" run - i.e. exactly
# where Word leaves _GoBack after you type a correction. The same
# bookmark is removed from its old location (after the "Mediante el
# menu..." paragraph) because Word only ever keeps one "_GoBack"
# bookmark at a time - re-adding it here automatically drops the old
# one.

$d = $word.ActiveDocument

# --- locate the text to fix -------------------------------------------------
$nameRange = $d.Content
$found = $nameRange.Find.Execute("Andres Fenrnao Galvis")
if (-not $found) {
    throw "Could not find 'Andres Fenrnao Galvis' in the document"
}
$nameStart = $nameRange.Start
$nameEnd = $nameRange.End

# Locate the id number that immediately follows the name (it lives in
# its own run, separated from the name by a single space run) so the
# existing run boundary there can be restored after the edit.
$idRange = $d.Range($nameEnd, $nameEnd + 20)
$idRange.Find.Execute("201632930") | Out-Null
$idStart = $idRange.Start

# --- apply the text fix ------------------------------------------------------
# Split point: right before the final "o" of "Fenrnao".
$fenrnaoOffset = $nameStart + ("Andres Fenrna").Length

$insPoint = $d.Range($fenrnaoOffset, $fenrnaoOffset)
$insPoint.InsertBefore("nd")

# Every position from the insertion point onward shifted right by 2.
$ndEnd = $fenrnaoOffset + 2
$newNameEnd = $nameEnd + 2
$newIdStart = $idStart + 2

# --- restore the run layout --------------------------------------------------
# Inserting text merges same-formatted neighbouring runs in the
# paragraph into one; toggling a character property back to its own
# value forces the engine to re-split the paragraph at that boundary
# without changing the visible formatting, which is how we reproduce
# the run boundaries from the target OOXML:
#   "Andres Fenrna" | "nd" | "o Galvis" | " " | "201632930"
foreach ($boundary in @($fenrnaoOffset, $ndEnd, $newNameEnd, $newIdStart)) {
    $splitRange = $d.Range($nameStart, $boundary)
    $splitRange.Bold = 1
    $splitRange.Bold = 0
}

# --- move the _GoBack bookmark ----------------------------------------------
# Placing it here both creates it at the right spot (between the new
# "nd" run and "o Galvis") and removes the old occurrence elsewhere in
# the document, matching the diff exactly.
$bmRange = $d.Range($ndEnd, $ndEnd)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
